# Apply the "data from jill, augmented by ian" edit to the Dragon Survey
# workbook:
#   1. Ten "Last rescue date" cells in column J were stored as text
#      (shared strings like "23/6/1634") with an obviously wrong century.
#      Ian corrected the century (+100 years) and turned them into real
#      Excel date serial numbers, matching the existing dd/mm/yyyy date
#      format already applied to the rest of column J. This also makes
#      the corresponding shared strings unused, so they disappear from
#      sharedStrings.xml and every other shared-string reference shifts
#      down automatically when the workbook is saved.
#   2. Column J was widened (23.33 -> 38 characters).
#   3. The active selection / scroll position on the sheet moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the ten mis-typed "Last rescue date" cells ---------------------
# row -> corrected Excel date serial (dd/mm/yyyy display, matches existing
# number format already applied to column J)
$fixes = @{
    6  = 12593   # 23/6/1634 -> 23/06/1934
    9  = 32999   # 6/5/1890  -> 06/05/1990
    10 = 32178   # 5/2/1888  -> 05/02/1988
    13 = 17015   # 1/8/1876  -> 01/08/1946
    22 = 18401   # 18/5/1850 -> 18/05/1950
    24 = 23074   # 4/3/1763  -> 04/03/1963
    26 = 20547   # 2/4/1856  -> 02/04/1956
    35 = 36406   # 3/9/1899  -> 03/09/1999
    36 = 27881   # 1/5/1876  -> 01/05/1976
    40 = 25906   # 4/12/1870 -> 04/12/1970
}

foreach ($row in $fixes.Keys) {
    $ws.Range("J$row").Value = $fixes[$row]
}

# --- 2. Widen column J -------------------------------------------------
$ws.Columns("J").ColumnWidth = 37.17

# --- 3. Update the sheet selection / scroll position --------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 8
$ws.Range("N26").Select()
